# Add description fields to Headers and all HDUs.
# Reads the HDU -> (Header description, Table description) mapping and
# writes it into a brand-new "HDU Descriptions" worksheet placed right
# after "L2A_product_definition".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the new sheet right after L2A_product_definition ---------------
$hdu = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$hdu.Name = "HDU Descriptions"

# Column widths (characters)
$hdu.Columns.Item(1).ColumnWidth = 20.67
$hdu.Columns.Item(2).ColumnWidth = 103.5
$hdu.Columns.Item(3).ColumnWidth = 35.83

# --- header row --------------------------------------------------------
$hdu.Range("B1").Value = "Header"
$hdu.Range("C1").Value = "Table"

# --- HDU description rows ------------------------------------------------
# Values are entered in the same order the original authors typed them in
# (not simple row order) so that shared-string indices line up with the
# source workbook.
$hdu.Range("B3").Value = "Configuration Header"
$hdu.Range("C3").Value = "Information on instrument configuration"

$hdu.Range("B6").Value = "Rawcounts and calibrated UVIS data"
$hdu.Range("C6").Value = "Rawcounts and calibrated UVIS data"

$hdu.Range("B11").Value = "2D arrays for each detector image in the FUV channel"
$hdu.Range("C11").Value = "2D arrays for each detector image in the FUV channel"

$hdu.Range("B10").Value = "This is the FITS header for the FOV_Geom HDU"
$hdu.Range("C10").Value = "Field Of View Geometry"

$hdu.Range("B2").Value = "This is the FITS header for the KERNELS HDU, which is list of all spice kernels used to compute geometry for this obsrevation"
$hdu.Range("C2").Value = "List of all kernels used to create the geomtrey"

$hdu.Range("B8").Value = "This is the FITS header for the SC_Geom HDU"
$hdu.Range("C8").Value = "Spacecraft Geometry"

$hdu.Range("B9").Value = "This is the FITS header for the Target_Geom HDU"
$hdu.Range("C9").Value = "Target Geometry"

$hdu.Range("B4").Value = "Time conversion from ET to UTC"
$hdu.Range("C4").Value = "Time conversion from ET to UTC"

$hdu.Range("B5").Value = "This is the FITS header for the Wavelength HDU"
$hdu.Range("C5").Value = "Wavelength calibration for EUV or FUV channel"

$hdu.Range("B7").Value = "Calibration factor"
$hdu.Range("C7").Value = "UVIS data calibration matrix"

$hdu.Range("B12").Value = "2D arrays for each detector image in the EUV channel"
$hdu.Range("C12").Value = "2D arrays for each detector image in the EUV channel"

# --- HDU NAME column (reuses existing shared strings) -----------------
$hdu.Range("A1").Value = "HDU NAME"
$hdu.Range("A2").Value = "KERNELS"
$hdu.Range("A3").Value = "CONFIG"
$hdu.Range("A4").Value = "TIME"
$hdu.Range("A5").Value = "WAVELENGTH"
$hdu.Range("A6").Value = "DATA"
$hdu.Range("A7").Value = "CAL"
$hdu.Range("A8").Value = "SC_GEOM"
$hdu.Range("A9").Value = "TARGET_GEOM"
$hdu.Range("A10").Value = "FOV_GEOM"
$hdu.Range("A11").Value = "DETECTOR_IMG_FUV"
$hdu.Range("A12").Value = "DETECTOR_IMG_EUV"

# --- formatting: bold header row ----------------------------------------
$hdu.Range("A1:C1").Font.Bold = $true

# --- tidy up the original sheet's lingering selection -------------------
$ws1.Range("B52").Select()

# --- view state: select B30 on the new sheet and make it the active tab -
$hdu.Range("B30").Select()
$hdu.Activate()
